$wb = $excel.ActiveWorkbook

# --- Sheet "About" ---
$about = $wb.Worksheets.Item("About")

# Source block
$about.Range("B3").Value = "European Commission"
$about.Range("B4").Value = 2016
$about.Range("B5").Value = "EU Reference Scenario 2016: Energy, transport and GHG emissions Trends to 2050 "
$about.Range("B6").Value = "https://ec.europa.eu/energy/sites/ener/files/documents/20160713%20draft_publication_REF2016_v13.pdf"
$about.Range("B7").Value = "Page 47 (see Notes below)"

# Notes block
$about.Range("A11").Value = "Vehicle buyer discount rates vary tremendously by study."
$about.Range("A20").Value = "We have chosen to use the discount rates applied in the PRIMES Reference Scenario 2016"
$about.Range("A21").Value = "These discount rates are differentiated by transport type (see below):"
$about.Range("A22").Value = "'- Public transport (road and conventional rail) 7.5%;"
$about.Range("A23").Value = "'- Public transport (advanced technologies, e.g. high speed rail) 8.5%;"
$about.Range("A24").Value = "'- Business transport sectors (aviation, trucks, maritime) 9.5% ;"
$about.Range("A25").Value = "'- Private cars and two-wheelers 11%."
$about.Range("A26").Value = ""
$about.Range("A27").Value = ""

# Hyperlink on B6 (added after the quote-prefixed cells so style indices line up)
$about.Hyperlinks.Add($about.Range("B6"), "https://ec.europa.eu/energy/sites/ener/files/documents/20160713%20draft_publication_REF2016_v13.pdf") | Out-Null

# --- Sheet "VBDR" ---
$vbdr = $wb.Worksheets.Item("VBDR")
$vbdr.Range("B2").Value = 0.11
$vbdr.Range("B3").Value = 0.095
$vbdr.Range("B4").Value = 0.095
$vbdr.Range("B5").Value = 0.075
$vbdr.Range("B6").Value = 0.095
$vbdr.Range("B7").Value = 0.11

# Restore the on-screen selections recorded for each sheet. VBDR's
# selection is set first, and "About" last, so that "About" ends up
# as the active (tabSelected) sheet, matching the saved workbook state.
$vbdr.Range("D5").Select() | Out-Null
$about.Range("M11").Select() | Out-Null
